$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.046810264341325
$ws.Range("D2").Value = 1.055234419612457
$ws.Range("E2").Value = 1.060410840985502
$ws.Range("F2").Value = 1.067775337790392
$ws.Range("I2").Value = 1.046217627763845
$ws.Range("J2").Value = 1.051862675701507
$ws.Range("K2").Value = 1.057975437745309
$ws.Range("L2").Value = 1.063137694830238
$ws.Range("M2").Value = 1.070482296160417
$ws.Range("N2").Value = 1.021136860290057
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.047641194868147
$ws.Range("D3").Value = 1.055895854492617
$ws.Range("E3").Value = 1.061195055456298
$ws.Range("F3").Value = 1.068565549263242
$ws.Range("I3").Value = 1.04641814162123
$ws.Range("J3").Value = 1.052342408947657
$ws.Range("K3").Value = 1.058450661547886
$ws.Range("L3").Value = 1.063736405663498
$ws.Range("M3").Value = 1.071088421744308
$ws.Range("N3").Value = 1.021298032604335
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.048179646282333
$ws.Range("D4").Value = 1.056324534269558
$ws.Range("E4").Value = 1.06170358757896
$ws.Range("F4").Value = 1.069077929518894
$ws.Range("I4").Value = 1.046547126965222
$ws.Range("J4").Value = 1.052652923715695
$ws.Range("K4").Value = 1.058758164043797
$ws.Range("L4").Value = 1.064124242425581
$ws.Range("M4").Value = 1.071481024376149
$ws.Range("N4").Value = 1.021402309999712
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.048406197120103
$ws.Range("D5").Value = 1.056504913880519
$ws.Range("E5").Value = 1.06191763412661
$ws.Range("F5").Value = 1.069293586037498
$ws.Range("I5").Value = 1.046601169564096
$ws.Range("J5").Value = 1.052783485594828
$ws.Range("K5").Value = 1.05888743690889
$ws.Range("L5").Value = 1.064287390517971
$ws.Range("M5").Value = 1.07164616833602
$ws.Range("N5").Value = 1.021446144875482
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.048444246850061
$ws.Range("D6").Value = 1.056535209893092
$ws.Range("E6").Value = 1.061953588651137
$ws.Range("F6").Value = 1.069329810425466
$ws.Range("I6").Value = 1.046610232810166
$ws.Range("J6").Value = 1.052805408728957
$ws.Range("K6").Value = 1.058909142281513
$ws.Range("L6").Value = 1.064314789706481
$ws.Range("M6").Value = 1.071673902186216
$ws.Range("N6").Value = 1.021453504733541
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.048182672737623
$ws.Range("D7").Value = 1.056326943874074
$ws.Range("E7").Value = 1.061706446661701
$ws.Range("F7").Value = 1.069080810144788
$ws.Range("I7").Value = 1.046547849804463
$ws.Range("J7").Value = 1.052654668207226
$ws.Range("K7").Value = 1.058759891400559
$ws.Range("L7").Value = 1.064126422021514
$ws.Range("M7").Value = 1.071483230671296
$ws.Range("N7").Value = 1.021402895737223
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.047090917840839
$ws.Range("D8").Value = 1.055457811211971
$ws.Range("E8").Value = 1.060675642436109
$ws.Range("F8").Value = 1.068042172622862
$ws.Range("I8").Value = 1.046285549371484
$ws.Range("J8").Value = 1.052024783264579
$ws.Range("K8").Value = 1.058136040819933
$ws.Range("L8").Value = 1.063339941780928
$ws.Range("M8").Value = 1.070687055549951
$ws.Range("N8").Value = 1.021191331385491
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.04517318807228
$ws.Range("D9").Value = 1.053931634529235
$ws.Range("E9").Value = 1.058867692835685
$ws.Range("F9").Value = 1.066220174134707
$ws.Range("I9").Value = 1.045817551770982
$ws.Range("J9").Value = 1.050915634516733
$ws.Range("K9").Value = 1.057036805560792
$ws.Range("L9").Value = 1.061957434404395
$ws.Range("M9").Value = 1.069287221942961
$ws.Range("N9").Value = 1.020818459161775
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.043898897077914
$ws.Range("D10").Value = 1.052917890235788
$ws.Range("E10").Value = 1.057668200433336
$ws.Range("F10").Value = 1.065011153352059
$ws.Range("I10").Value = 1.045501706009109
$ws.Range("J10").Value = 1.050176812578178
$ws.Range("K10").Value = 1.056304112538769
$ws.Range("L10").Value = 1.061038126438448
$ws.Range("M10").Value = 1.068356203076759
$ws.Range("N10").Value = 1.020569861274243
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.04334813260165
$ws.Range("D11").Value = 1.052479830794264
$ws.Range("E11").Value = 1.057150209799574
$ws.Range("F11").Value = 1.064488998864144
$ws.Range("I11").Value = 1.045364038488716
$ws.Range("J11").Value = 1.049857056278582
$ws.Range("K11").Value = 1.055986897093258
$ws.Range("L11").Value = 1.060640637082423
$ws.Range("M11").Value = 1.067953605450142
$ws.Range("N11").Value = 1.020462218069341
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.043143707998282
$ws.Range("D12").Value = 1.052317252944367
$ws.Range("E12").Value = 1.056958017040127
$ws.Range("F12").Value = 1.064295253800764
$ws.Range("I12").Value = 1.045312767567444
$ws.Range("J12").Value = 1.049738309729385
$ws.Range("K12").Value = 1.055869077419937
$ws.Range("L12").Value = 1.060493080197917
$ws.Range("M12").Value = 1.067804145564075
$ws.Range("N12").Value = 1.020422235299028
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.043187550778623
$ws.Range("D13").Value = 1.052352120217533
$ws.Range("E13").Value = 1.05699923341165
$ws.Range("F13").Value = 1.064336803402273
$ws.Range("I13").Value = 1.045323771454053
$ws.Range("J13").Value = 1.049763780105966
$ws.Range("K13").Value = 1.055894349742648
$ws.Range("L13").Value = 1.060524727625325
$ws.Range("M13").Value = 1.067836201442646
$ws.Range("N13").Value = 1.02043081169872
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.043331231651249
$ws.Range("D14").Value = 1.052466389253919
$ws.Range("E14").Value = 1.057134318748912
$ws.Range("F14").Value = 1.064472979614811
$ws.Range("I14").Value = 1.045359803169362
$ws.Range("J14").Value = 1.04984724013025
$ws.Range("K14").Value = 1.055977157905354
$ws.Range("L14").Value = 1.060628438172169
$ws.Range("M14").Value = 1.067941249350353
$ws.Range("N14").Value = 1.020458913063664
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.04341977861935
$ws.Range("D15").Value = 1.052536812381241
$ws.Range("E15").Value = 1.057217577440853
$ws.Range("F15").Value = 1.064556909681983
$ws.Range("I15").Value = 1.045381985622875
$ws.Range("J15").Value = 1.049898665982962
$ws.Range("K15").Value = 1.05602817988585
$ws.Range("L15").Value = 1.060692349421355
$ws.Range("M15").Value = 1.068005983856201
$ws.Range("N15").Value = 1.020476227351251
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.043935470977368
$ws.Range("D16").Value = 1.052946981894589
$ws.Range("E16").Value = 1.057702607385839
$ws.Range("F16").Value = 1.065045835849136
$ws.Range("I16").Value = 1.045510823564197
$ws.Range("J16").Value = 1.050198037206385
$ws.Range("K16").Value = 1.056325166163701
$ws.Range("L16").Value = 1.061064518797918
$ws.Range("M16").Value = 1.068382933672392
$ws.Range("N16").Value = 1.020577005270638
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.044259223372123
$ws.Range("D17").Value = 1.053204512323495
$ws.Range("E17").Value = 1.058007229464485
$ws.Range("F17").Value = 1.06535289178337
$ws.Range("I17").Value = 1.04559139866128
$ws.Range("J17").Value = 1.050385868259798
$ws.Range("K17").Value = 1.056511470996308
$ws.Range("L17").Value = 1.061298126485084
$ws.Range("M17").Value = 1.068619530034187
$ws.Range("N17").Value = 1.020640221310969
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.044448160181554
$ws.Range("D18").Value = 1.053354812030603
$ws.Range("E18").Value = 1.058185045068221
$ws.Range("F18").Value = 1.065532123344086
$ws.Range("I18").Value = 1.045638309455895
$ws.Range("J18").Value = 1.050495442134168
$ws.Range("K18").Value = 1.056620143743823
$ws.Range("L18").Value = 1.061434441462995
$ws.Range("M18").Value = 1.068757584645923
$ws.Range("N18").Value = 1.020677094251003
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.044512599236099
$ws.Range("D19").Value = 1.053406074970582
$ws.Range("E19").Value = 1.058245698378294
$ws.Range("F19").Value = 1.065593258835023
$ws.Range("I19").Value = 1.045654289994079
$ws.Range("J19").Value = 1.050532806530033
$ws.Range("K19").Value = 1.056657199026459
$ws.Range("L19").Value = 1.06148093072216
$ws.Range("M19").Value = 1.06880466646143
$ws.Range("N19").Value = 1.020689666970956
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.04422447770461
$ws.Range("D20").Value = 1.053176872777921
$ws.Range("E20").Value = 1.057974532450894
$ws.Range("F20").Value = 1.06531993401936
$ws.Range("I20").Value = 1.045582762739804
$ws.Range("J20").Value = 1.050365714183435
$ws.Range("K20").Value = 1.056491481796781
$ws.Range("L20").Value = 1.061273056831611
$ws.Range("M20").Value = 1.068594140097183
$ws.Range("N20").Value = 1.020633438814787
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.043288916953923
$ws.Range("D21").Value = 1.052432736075257
$ws.Range("E21").Value = 1.057094533614456
$ws.Range("F21").Value = 1.064432873400733
$ws.Range("I21").Value = 1.045349196449759
$ws.Range("J21").Value = 1.049822662527266
$ws.Range("K21").Value = 1.055952772725627
$ws.Range("L21").Value = 1.060597895548495
$ws.Range("M21").Value = 1.067910313060366
$ws.Range("N21").Value = 1.020450637888433
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.042701583783334
$ws.Range("D22").Value = 1.051965660995579
$ws.Range("E22").Value = 1.056542471155918
$ws.Range("F22").Value = 1.063876338177391
$ws.Range("I22").Value = 1.045201562938183
$ws.Range("J22").Value = 1.049481370745088
$ws.Range("K22").Value = 1.055614113493579
$ws.Range("L22").Value = 1.060173906903089
$ws.Range("M22").Value = 1.067480843755208
$ws.Range("N22").Value = 1.020335708043978
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.04301285514595
$ws.Range("D23").Value = 1.052213190499928
$ws.Range("E23").Value = 1.056835012844947
$ws.Range("F23").Value = 1.064171254051772
$ws.Range("I23").Value = 1.045279900037484
$ws.Range("J23").Value = 1.049662281647112
$ws.Range("K23").Value = 1.055793638079282
$ws.Range("L23").Value = 1.060398622214125
$ws.Range("M23").Value = 1.067708467498209
$ws.Range("N23").Value = 1.020396633946099
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.044240177475076
$ws.Range("D24").Value = 1.053189361631624
$ws.Range("E24").Value = 1.057989306407988
$ws.Range("F24").Value = 1.06533482580969
$ws.Range("I24").Value = 1.045586665210822
$ws.Range("J24").Value = 1.050374820896239
$ws.Range("K24").Value = 1.056500514042554
$ws.Range("L24").Value = 1.061284384557434
$ws.Range("M24").Value = 1.068605612556507
$ws.Range("N24").Value = 1.020636503532641
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.045668235501829
$ws.Range("D25").Value = 1.054325542654441
$ws.Range("E25").Value = 1.059334076336821
$ws.Range("F25").Value = 1.066690218440396
$ws.Range("I25").Value = 1.045939221867989
$ws.Range("J25").Value = 1.051202274304454
$ws.Range("K25").Value = 1.057320967375311
$ws.Range("L25").Value = 1.062314436194957
$ws.Range("M25").Value = 1.069648731416974
$ws.Range("N25").Value = 1.020914860645799
